{"js": "// Replace the two-digit multiplication expressions in the document's table\n// with the new values from the commit's diff. Each old expression is\n// unique within the document, so a simple search/replace per pair is safe\n// and preserves run formatting (font, size) since we replace only the text\n// inside the existing run via Range.insertText(..., \"Replace\").\nconst replacements = [\n  [\"82\u00d738=\", \"17\u00d749=\"],\n  [\"76\u00d767=\", \"11\u00d796=\"],\n  [\"68\u00d739=\", \"62\u00d732=\"],\n  [\"69\u00d790=\", \"94\u00d741=\"],\n  [\"59\u00d717=\", \"67\u00d761=\"],\n  [\"61\u00d765=\", \"87\u00d728=\"],\n  [\"27\u00d789=\", \"83\u00d711=\"],\n  [\"19\u00d797=\", \"42\u00d716=\"],\n  [\"82\u00d751=\", \"67\u00d718=\"],\n  [\"47\u00d761=\", \"69\u00d724=\"],\n  [\"39\u00d766=\", \"83\u00d763=\"],\n  [\"43\u00d753=\", \"24\u00d741=\"],\n  [\"98\u00d750=\", \"47\u00d793=\"],\n  [\"24\u00d720=\", \"51\u00d774=\"],\n  [\"91\u00d760=\", \"28\u00d750=\"],\n  [\"26\u00d775=\", \"54\u00d753=\"],\n  [\"23\u00d736=\", \"58\u00d746=\"],\n  [\"13\u00d775=\", \"34\u00d714=\"],\n  [\"69\u00d789=\", \"96\u00d723=\"],\n  [\"72\u00d764=\", \"96\u00d764=\"],\n  [\"43\u00d760=\", \"95\u00d760=\"],\n  [\"84\u00d749=\", \"78\u00d725=\"],\n  [\"82\u00d752=\", \"92\u00d722=\"],\n  [\"13\u00d733=\", \"36\u00d790=\"],\n  [\"52\u00d775=\", \"27\u00d780=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication expressions in the document's table\n# with the new values from the commit's diff. Each old expression is\n# unique within the document, so Find/Replace (wdReplaceAll = 2) per pair\n# safely updates only the matching run text and preserves formatting.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"82\u00d738=\", \"17\u00d749=\"),\n  @(\"76\u00d767=\", \"11\u00d796=\"),\n  @(\"68\u00d739=\", \"62\u00d732=\"),\n  @(\"69\u00d790=\", \"94\u00d741=\"),\n  @(\"59\u00d717=\", \"67\u00d761=\"),\n  @(\"61\u00d765=\", \"87\u00d728=\"),\n  @(\"27\u00d789=\", \"83\u00d711=\"),\n  @(\"19\u00d797=\", \"42\u00d716=\"),\n  @(\"82\u00d751=\", \"67\u00d718=\"),\n  @(\"47\u00d761=\", \"69\u00d724=\"),\n  @(\"39\u00d766=\", \"83\u00d763=\"),\n  @(\"43\u00d753=\", \"24\u00d741=\"),\n  @(\"98\u00d750=\", \"47\u00d793=\"),\n  @(\"24\u00d720=\", \"51\u00d774=\"),\n  @(\"91\u00d760=\", \"28\u00d750=\"),\n  @(\"26\u00d775=\", \"54\u00d753=\"),\n  @(\"23\u00d736=\", \"58\u00d746=\"),\n  @(\"13\u00d775=\", \"34\u00d714=\"),\n  @(\"69\u00d789=\", \"96\u00d723=\"),\n  @(\"72\u00d764=\", \"96\u00d764=\"),\n  @(\"43\u00d760=\", \"95\u00d760=\"),\n  @(\"84\u00d749=\", \"78\u00d725=\"),\n  @(\"82\u00d752=\", \"92\u00d722=\"),\n  @(\"13\u00d733=\", \"36\u00d790=\"),\n  @(\"52\u00d775=\", \"27\u00d780=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
